# 7RMartSupermarket/testdata.xlsx - minor changes and new Testng file creation
#
# Content edits (3 cells get new text; Excel's save path garbage-collects the
# now-unreferenced shared strings and appends the replacements at the end of
# the shared-string table, which is what the target diff shows):
#   - ManageNewsPage!B2       : "Hello, Welcome to Coding!!" -> "Welcome to Coding!! JS"
#   - ManageFooterTextPage!A2 : "Fno: 12B, Skyline, Kerala"  -> "FNmbr: 1207B, Skyline Apart, Kerala"
#   - ManageCategory!A2       : "Pizza XL"                   -> "Pizza XS"
#
# Plus cursor/selection + active-tab bookkeeping: the active cell moves on
# several sheets, and the active/selected tab moves from LoginPage to
# ManageContactPage (last sheet selected becomes tabSelected + workbook
# activeTab).

$wb = $excel.ActiveWorkbook

$wsNews = $wb.Worksheets.Item("ManageNewsPage")
$wsNews.Range("B2").Value = "Welcome to Coding!! JS"

$wsFooter = $wb.Worksheets.Item("ManageFooterTextPage")
$wsFooter.Range("A2").Value = "FNmbr: 1207B, Skyline Apart, Kerala"

$wsCategory = $wb.Worksheets.Item("ManageCategory")
$wsCategory.Range("A2").Value = "Pizza XS"

# Move each sheet's active-cell selection to match the edited workbook, in
# sheet order, finishing on ManageContactPage so it becomes the active /
# tabSelected sheet (matching activeTab="5" on the workbook view and the
# tabSelected flag moving off LoginPage and onto ManageContactPage).
$wsNews.Range("F4").Select()
$wsFooter.Range("D7").Select()
$wsCategory.Range("B5").Select()

$wsAdmin = $wb.Worksheets.Item("AdminUsers")
$wsAdmin.Range("A3").Select()

$wsContact = $wb.Worksheets.Item("ManageContactPage")
$wsContact.Range("A6").Select()
